# Replace the M2Doc `m:'doc.html'.fromHTMLURI()` Word FIELD (fldChar/instrText)
# in the second paragraph with literal `{ ... }` text runs, as produced by the
# new TokenIteratorFieldRewriterSplit parser (fields are no longer used to
# store M2Doc query tokens; the query text is written directly as w:t runs).

$d = $word.ActiveDocument

# Locate the paragraph that holds the M2Doc field (robust to index drift).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing the M2Doc field"
}

# Grab the field so we can sanity check / reuse its code text.
$field = $target.Range.Fields(1)
$code = $field.Code.Text   # " m:'doc.html'.fromHTMLURI() "

# Build the replacement OOXML: the field's begin mark + leading space become a
# single literal "{" run, the quoted-string tokens become individual w:t runs
# (so the existing _GoBack bookmark keeps sitting between "doc.html" and the
# rest of the literal), and the trailing space + end mark become a single "}"
# run.
$openBrace  = "{"
$closeBrace = "}"
$tok1 = "m"
$tok2 = ":"
$tok3 = "'"
$tok4 = "doc.html"
$tok5 = "'.fromHTMLURI()"

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>$openBrace</w:t></w:r>
<w:r><w:t>$tok1</w:t></w:r>
<w:r><w:t>$tok2</w:t></w:r>
<w:r><w:t>$tok3</w:t></w:r>
<w:r><w:t>$tok4</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:t>$tok5</w:t></w:r>
<w:r><w:t xml:space="preserve">$closeBrace</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($xml)

Write-Output "Rewrote field code [$code] into literal braces text runs."
